$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '64.897.59'
Set-TextValue $ws 'E2' '  -1.99%  '
Set-TextValue $ws 'D3' '3.236.55'
Set-TextValue $ws 'E3' '  -1.29%  '
Set-TextValue $ws 'E4' '  -0.03%  '
Set-TextValue $ws 'D5' '577.93'
Set-TextValue $ws 'E5' '  -0.38%  '
Set-TextValue $ws 'D6' '172.81'
Set-TextValue $ws 'E6' '  -3.55%  '
Set-TextValue $ws 'E7' '  +0.29%  '
Set-TextValue $ws 'E8' '  -0.01%  '
Set-TextValue $ws 'D9' '3.234.33'
Set-TextValue $ws 'E9' '  -1.33%  '
Set-TextValue $ws 'E10' '  -2.72%  '
Set-TextValue $ws 'E11' '  +0.66%  '
Set-TextValue $ws 'E12' '  -3.10%  '
Set-TextValue $ws 'D13' '3.796.45'
Set-TextValue $ws 'E13' '  -1.43%  '
Set-TextValue $ws 'E14' '  -3.18%  '
Set-TextValue $ws 'D15' '64.945.66'
Set-TextValue $ws 'E15' '  -1.96%  '
Set-TextValue $ws 'D16' '25.79'
Set-TextValue $ws 'E16' '  -2.08%  '
Set-TextValue $ws 'B17' 'WrappedEther'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D17' '3.261.87'
Set-TextValue $ws 'E17' '  -0.58%  '
Set-TextValue $ws 'B18' 'ShibaInu'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws 'D18' '0.0000159'
Set-TextValue $ws 'E18' '  -2.87%  '
Set-TextValue $ws 'D19' '418.90'
Set-TextValue $ws 'E19' '  -3.81%  '
Set-TextValue $ws 'E20' '  -2.26%  '
Set-TextValue $ws 'D21' '12.86'
Set-TextValue $ws 'E21' '  -2.40%  '
Set-TextValue $ws 'D22' '7.21'
Set-TextValue $ws 'E22' '  -2.50%  '
Set-TextValue $ws 'E23' '  -0.03%  '
Set-TextValue $ws 'D24' '70.95'
Set-TextValue $ws 'E24' '  -1.09%  '
Set-TextValue $ws 'D25' '5.67'
Set-TextValue $ws 'E25' '  -0.09%  '
Set-TextValue $ws 'E26' '  +4.14%  '
Set-TextValue $ws 'E27' '  -1.95%  '
Set-TextValue $ws 'D28' '0.0000112'
Set-TextValue $ws 'E28' '  -1.14%  '
Set-TextValue $ws 'D29' '9.11'
Set-TextValue $ws 'E29' '  +2.97%  '
Set-TextValue $ws 'D30' '0.999'
Set-TextValue $ws 'E31' '  -3.78%  '
Set-TextValue $ws 'D32' '21.87'
Set-TextValue $ws 'E32' '  -1.93%  '
Set-TextValue $ws 'E33' '  +0.04%  '
Set-TextValue $ws 'D34' '5.01'
Set-TextValue $ws 'E34' '  -3.58%  '
Set-TextValue $ws 'E35' '  -2.28%  '
Set-TextValue $ws 'E36' '  -1.81%  '
Set-TextValue $ws 'D37' '157.52'
Set-TextValue $ws 'E37' '  -0.38%  '
Set-TextValue $ws 'E38' '  -1.89%  '
Set-TextValue $ws 'D39' '2.823.85'
Set-TextValue $ws 'E39' '  +1.74%  '
Set-TextValue $ws 'E40' '  -2.89%  '
Set-TextValue $ws 'D41' '25.49'
Set-TextValue $ws 'E41' '  -3.92%  '
Set-TextValue $ws 'E42' '  -1.80%  '
Set-TextValue $ws 'D43' '39.60'
Set-TextValue $ws 'E43' '  -1.64%  '
Set-TextValue $ws 'D44' '0.725'
Set-TextValue $ws 'E44' '  -6.21%  '
Set-TextValue $ws 'D45' '5.77'
Set-TextValue $ws 'E45' '  -4.29%  '
Set-TextValue $ws 'D46' '0.0630'
Set-TextValue $ws 'E46' '  -4.38%  '
Set-TextValue $ws 'D47' '304.12'
Set-TextValue $ws 'E47' '  -5.53%  '
Set-TextValue $ws 'D48' '2.16'
Set-TextValue $ws 'E48' '  -5.28%  '
Set-TextValue $ws 'D49' '22.15'
Set-TextValue $ws 'E49' '  -4.56%  '
Set-TextValue $ws 'E50' '  -1.00%  '
Set-TextValue $ws 'E51' '  -1.11%  '
